# Weekly roll of the "Acelga" price series: a new observation is inserted
# at row 72 and every subsequent observation (rows 72..201) shifts down one
# row, with the row that used to be last (201) becoming the new last row
# (202). Columns A,B,C,E,F,G,H,N,O,Q,R are constant for every data row in
# this range, so only D (Fecha), I (Calidad), J (Volumen), K/L/M (precios),
# P (Precio $/Kg) actually need to move; row 202 is brand new so every
# column must be populated there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colD = 4   # Fecha
$colI = 9   # Calidad
$colJ = 10  # Volumen
$colK = 11  # Precio minimo
$colL = 12  # Precio maximo
$colM = 13  # Precio promedio ponderado
$colN = 14  # Unidad de comercializacion
$colP = 16  # Precio $/Kg

$firstShift = 72
$lastShiftOld = 201   # old last data row; becomes row 202 after the shift

# --- 1. Snapshot the "old" values for rows 71..201 before overwriting anything ---
$oldD = @{}
$oldI = @{}
$oldJ = @{}
$oldK = @{}
$oldL = @{}
$oldM = @{}
$oldP = @{}

for ($r = ($firstShift - 1); $r -le $lastShiftOld; $r++) {
    $oldD[$r] = $ws.Cells.Item($r, $colD).Value2
    $oldI[$r] = $ws.Cells.Item($r, $colI).Value2
    $oldJ[$r] = $ws.Cells.Item($r, $colJ).Value2
    $oldK[$r] = $ws.Cells.Item($r, $colK).Value2
    $oldL[$r] = $ws.Cells.Item($r, $colL).Value2
    $oldM[$r] = $ws.Cells.Item($r, $colM).Value2
    $oldP[$r] = $ws.Cells.Item($r, $colP).Value2
}

# --- 2. Row 202 is new: duplicate row 201 in full (all columns A..R) ---
$ws.Range("A201:R201").Copy()
$ws.Range("A202:R202").PasteSpecial()

# --- 3. Shift rows 72..201 down by one: row r gets what row (r-1) used to hold ---
for ($r = $lastShiftOld; $r -ge ($firstShift + 1); $r--) {
    $src = $r - 1
    $ws.Cells.Item($r, $colD).Value = $oldD[$src]
    $ws.Cells.Item($r, $colJ).Value = $oldJ[$src]
    $ws.Cells.Item($r, $colK).Value = $oldK[$src]
    $ws.Cells.Item($r, $colL).Value = $oldL[$src]
    $ws.Cells.Item($r, $colM).Value = $oldM[$src]
    $ws.Cells.Item($r, $colP).Value = $oldP[$src]
    if ($oldI[$src] -ne $oldI[$r]) {
        $ws.Cells.Item($r, $colI).Value = $oldI[$src]
    }
}

# --- 4. Row 72 becomes the brand-new observation (not shifted from anything) ---
$ws.Cells.Item($firstShift, $colD).Value = 44469
$ws.Cells.Item($firstShift, $colJ).Value = 141
$ws.Cells.Item($firstShift, $colK).Value = 2000
$ws.Cells.Item($firstShift, $colL).Value = 2200
$ws.Cells.Item($firstShift, $colM).Value = 2184
$ws.Cells.Item($firstShift, $colP).Value = 364

# --- 5. Make sure every "Fecha" cell keeps its date number format ---
$ws.Range("D72:D202").NumberFormat = "YYYY-MM-DD HH:MM:SS"
